$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)

# Add paragraph border (top/left/bottom/right) with 5pt space-from-text on each side,
# matching <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$b = $p.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# Replace the paragraph's text (excluding the trailing paragraph mark) with the new
# ID marker, dropping the separate trailing-space run entirely.
$r = $p.Range
$r.End = $r.End - 1
$r.Text = "**ID__AFFARS_AFFARS_PGI_5308__ID**"
